$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 95, pushing the existing rows 95-120
# down to 97-122 (dimension grows from A1:T120 to A1:T122).
$ws.Rows("95:96").Insert()

# --- New row 95 ---
$ws.Range("A95").Value = 3
$ws.Range("B95").Value = "Femacal de La Calera"
$ws.Range("C95").Value = "Coquimbo"
$ws.Range("D95").Value = 44637
$ws.Range("E95").Value = 5
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100107
$ws.Range("H95").Value = "Otros"
$ws.Range("I95").Value = 100107011
$ws.Range("J95").Value = "Tuna"
$ws.Range("K95").Value = "Sin especificar"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 60
$ws.Range("N95").Value = 16000
$ws.Range("O95").Value = 16000
$ws.Range("P95").Value = 16000
$ws.Range("Q95").Value = "`$/caja 16 kilos"
$ws.Range("R95").Value = "Cabildo"
$ws.Range("S95").Value = 1000
$ws.Range("T95").Value = 16

# --- New row 96 ---
$ws.Range("A96").Value = 3
$ws.Range("B96").Value = "Femacal de La Calera"
$ws.Range("C96").Value = "Coquimbo"
$ws.Range("D96").Value = 44637
$ws.Range("E96").Value = 5
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100107
$ws.Range("H96").Value = "Otros"
$ws.Range("I96").Value = 100107011
$ws.Range("J96").Value = "Tuna"
$ws.Range("K96").Value = "Sin especificar"
$ws.Range("L96").Value = "Segunda"
$ws.Range("M96").Value = 65
$ws.Range("N96").Value = 14000
$ws.Range("O96").Value = 14000
$ws.Range("P96").Value = 14000
$ws.Range("Q96").Value = "`$/caja 16 kilos"
$ws.Range("R96").Value = "Cabildo"
$ws.Range("S96").Value = 875
$ws.Range("T96").Value = 16
